# Add data for 2022-12-07
# 8 new crime incidents were recorded across the city on this date.
# This updates the aggregated "full year" rollups: the citywide totals
# sheet, the by-neighborhood summary sheet, and the per-neighborhood
# detail sheets for the neighborhoods/categories affected. Two of the
# affected neighborhoods (Rogers Park, Edgewater) had no prior Arson
# record for 2022, so a brand-new category row is inserted for each.

$wb = $excel.ActiveWorkbook

function Set-CellValue($ws, [string]$ref, $value) {
    $ws.Range($ref).Value = $value
}

# ---------------------------------------------------------------------
# Citywide Totals: Aggravated Battery +1, Arson +2, Robbery +2, Theft +3
# (Total +8)
# ---------------------------------------------------------------------
$wsCitywide = $wb.Worksheets.Item("Citywide Totals")
Set-CellValue $wsCitywide "I3" 194
Set-CellValue $wsCitywide "I4" 5
Set-CellValue $wsCitywide "I9" 499
Set-CellValue $wsCitywide "I10" 851
Set-CellValue $wsCitywide "I11" 1691

# ---------------------------------------------------------------------
# By Neighborhood: one new incident each in Chatham, Edgewater,
# Grand Crossing, Printers Row, Rogers Park, Uptown, Washington Park,
# Wicker Park (Total +8)
# ---------------------------------------------------------------------
$wsNeighborhood = $wb.Worksheets.Item("By Neighborhood")
Set-CellValue $wsNeighborhood "I19" 40
Set-CellValue $wsNeighborhood "I27" 23
Set-CellValue $wsNeighborhood "I36" 72
Set-CellValue $wsNeighborhood "I72" 16
Set-CellValue $wsNeighborhood "I76" 47
Set-CellValue $wsNeighborhood "I87" 23
Set-CellValue $wsNeighborhood "I89" 16
Set-CellValue $wsNeighborhood "I96" 15
Set-CellValue $wsNeighborhood "I99" 1691

# ---------------------------------------------------------------------
# Chatham: Theft +1 (Total +1)
# ---------------------------------------------------------------------
$wsChatham = $wb.Worksheets.Item("Chatham")
Set-CellValue $wsChatham "I9" 18
Set-CellValue $wsChatham "I10" 40

# ---------------------------------------------------------------------
# Grand Crossing: Robbery +1 (Total +1)
# ---------------------------------------------------------------------
$wsGrandCrossing = $wb.Worksheets.Item("Grand Crossing")
Set-CellValue $wsGrandCrossing "I7" 20
Set-CellValue $wsGrandCrossing "I9" 72

# ---------------------------------------------------------------------
# Washington Park: Aggravated Battery +1 (Total +1)
# ---------------------------------------------------------------------
$wsWashingtonPark = $wb.Worksheets.Item("Washington Park")
Set-CellValue $wsWashingtonPark "I3" 3
Set-CellValue $wsWashingtonPark "I6" 16

# ---------------------------------------------------------------------
# Uptown: Theft +1 (Total +1)
# ---------------------------------------------------------------------
$wsUptown = $wb.Worksheets.Item("Uptown")
Set-CellValue $wsUptown "I8" 13
Set-CellValue $wsUptown "I9" 23

# ---------------------------------------------------------------------
# Wicker Park: Robbery +1 (Total +1)
# ---------------------------------------------------------------------
$wsWickerPark = $wb.Worksheets.Item("Wicker Park")
Set-CellValue $wsWickerPark "I5" 4
Set-CellValue $wsWickerPark "I7" 15

# ---------------------------------------------------------------------
# Printers Row: Theft +1 (Total +1)
# ---------------------------------------------------------------------
$wsPrintersRow = $wb.Worksheets.Item("Printers Row")
Set-CellValue $wsPrintersRow "I5" 5
Set-CellValue $wsPrintersRow "I6" 16

# ---------------------------------------------------------------------
# Rogers Park: brand-new "Arson" category row (value 1) inserted at
# row 4, pushing existing category rows down by one. Total row +1.
# ---------------------------------------------------------------------
$wsRogersPark = $wb.Worksheets.Item("Rogers Park")
$wsRogersPark.Rows.Item(4).Insert()
$wsRogersPark.Range("A3").Copy()
$wsRogersPark.Range("A4").PasteSpecial(-4122)
$wsRogersPark.Cells.Item(4, 1).Value = "Arson"
$wsRogersPark.Cells.Item(4, 9).Value = 1
Set-CellValue $wsRogersPark "I10" 47

# ---------------------------------------------------------------------
# Edgewater: brand-new "Arson" category row (value 1) inserted at
# row 4, pushing existing category rows down by one. Total row +1.
# ---------------------------------------------------------------------
$wsEdgewater = $wb.Worksheets.Item("Edgewater")
$wsEdgewater.Rows.Item(4).Insert()
$wsEdgewater.Range("A3").Copy()
$wsEdgewater.Range("A4").PasteSpecial(-4122)
$wsEdgewater.Cells.Item(4, 1).Value = "Arson"
$wsEdgewater.Cells.Item(4, 9).Value = 1
Set-CellValue $wsEdgewater "I7" 23

$excel.CutCopyMode = $false
